# Updated cryptos list values (price & volume) per source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.833.43"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "1.985.50"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "239.63"
$ws.Range("E5").Value = "  -4.83%  "
$ws.Range("D6").Value = "0.602"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "53.86"
$ws.Range("E8").Value = "  -5.15%  "
$ws.Range("D9").Value = "0.371"
$ws.Range("E9").Value = "  -3.62%  "
$ws.Range("D10").Value = "57.48"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "0.0749"
$ws.Range("E11").Value = "  -4.53%  "
$ws.Range("D12").Value = "0.0976"
$ws.Range("E12").Value = "  -3.97%  "
$ws.Range("D13").Value = "2.278.39"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "14.03"
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").Value = "20.62"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").Value = "0.752"
$ws.Range("E16").Value = "  -7.69%  "
$ws.Range("D17").Value = "5.02"
$ws.Range("E17").Value = "  -6.04%  "
$ws.Range("D18").Value = "1.991.56"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "36.786.14"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").Value = "68.11"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").Value = "0.0₃0804"
$ws.Range("E21").Value = "  -5.07%  "
$ws.Range("D22").Value = "5.06"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").Value = "224.93"
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -7.12%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "162.11"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "8.62"
$ws.Range("E28").Value = "  -4.71%  "
$ws.Range("D29").Value = "19.07"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  -5.20%  "
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.40"
$ws.Range("E33").Value = "  -6.00%  "
$ws.Range("D34").Value = "0.0607"
$ws.Range("E34").Value = "  -8.12%  "
$ws.Range("D35").Value = "4.21"
$ws.Range("E35").Value = "  -7.12%  "
$ws.Range("E36").Value = "  -6.70%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").Value = "3.23"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("D40").Value = "5.27"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "1.427.84"
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0889"
$ws.Range("E43").Value = "  -7.67%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0203"
$ws.Range("E44").Value = "  -5.71%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "1.11"
$ws.Range("E45").Value = "  -5.20%  "
$ws.Range("D46").Value = "87.57"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("D48").Value = "14.83"
$ws.Range("E48").Value = "  -7.33%  "
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "6.71"
$ws.Range("E50").Value = "  -8.66%  "
$ws.Range("D51").Value = "2.170.41"
$ws.Range("E51").Value = "  -1.90%  "
